$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows 2-26, columns: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26)
$dCol = @(44839, 44841, 44845, 44826, 44749, 44824, 44518, 44756, 44757, 44827, 44830, 44825, 44769, 44812, 44819, 44525, 44776, 44838, 44813, 44508, 44767, 44755, 44811, 44771, 44837)
$jCol = @(80, 20, 20, 50, 50, 20, 50, 80, 30, 20, 25, 30, 50, 80, 100, 40, 80, 10, 20, 40, 50, 50, 30, 40, 80)
$kCol = @(16000, 16000, 16000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 12000, 20000, 20000, 20000, 20000, 8000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 20000, 16000)
$lCol = @(16000, 16000, 16000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 12000, 20000, 20000, 20000, 20000, 8000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 20000, 16000)
$mCol = @(16000, 16000, 16000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 12000, 20000, 20000, 20000, 20000, 8000, 20000, 20000, 20000, 10000, 20000, 20000, 20000, 20000, 16000)
$pCol = @(1067, 1067, 1067, 1333, 1333, 1333, 667, 1333, 1333, 1333, 800, 1333, 1333, 1333, 1333, 533, 1333, 1333, 1333, 667, 1333, 1333, 1333, 1333, 1067)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $dCol[$i]
    $ws.Cells.Item($r, 10).Value = $jCol[$i]
    $ws.Cells.Item($r, 11).Value = $kCol[$i]
    $ws.Cells.Item($r, 12).Value = $lCol[$i]
    $ws.Cells.Item($r, 13).Value = $mCol[$i]
    $ws.Cells.Item($r, 16).Value = $pCol[$i]
}

$wb.Save()
